{"js": "// The L/FNG Review template has several table-label cells whose paragraph\n// mark run properties (w:pPr/w:rPr) carry an accidentally duplicated\n// <w:b/> element (bold turned on twice). We touch each of those label\n// cells' paragraph formatting so the document model re-serializes the\n// (single, deduplicated) bold flag. We also update the \"Resolution for\n// Application to Proceed to the ALC\" label text to the new question\n// wording, per the commit \"Update L/FNG Review Form Labels\".\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load('items');\nawait context.sync();\n\n// [tableIndex, rowIndex, columnIndex] for every label cell that has the\n// duplicated bold flag in its paragraph mark run properties.\nconst boldLabelCells = [\n  [0, 0, 0], // Application ID:\n  [0, 1, 0], // Application Type:\n  [0, 2, 0], // Status:\n  [0, 3, 0], // Applicant:\n  [0, 4, 0], // Local/ First Nation Government:\n  [1, 0, 0], // Local Government File Number\n  [1, 1, 0], // First Name\n  [1, 2, 0], // Last Name\n  [1, 3, 0], // Position\n  [1, 4, 0], // Department\n  [1, 5, 0], // Phone Number\n  [1, 6, 0], // Email\n  [2, 0, 0], // Resolution for Application to Proceed to the ALC\n  [3, 0, 0], // Type (visible attachments table header)\n  [3, 0, 1], // File Name (visible attachments table header)\n  [4, 0, 0], // Type (no-data attachments table header)\n  [4, 0, 1], // File Name (no-data attachments table header)\n];\n\nconst firstParas = [];\nfor (const [tableIndex, rowIndex, columnIndex] of boldLabelCells) {\n  const cell = tables.items[tableIndex].getCell(rowIndex, columnIndex);\n  const para = cell.body.paragraphs.getFirst();\n  // Re-asserting bold forces the paragraph's formatting (w:pPr/w:rPr and\n  // w:r/w:rPr) to be rewritten without the stray duplicate <w:b/>.\n  para.font.bold = true;\n  firstParas.push(para);\n}\n\n// Grab the label text of the \"Resolution\" row so we can safely confirm\n// we're rewriting the right cell before changing it.\nconst resolutionPara = firstParas[12];\nresolutionPara.load('text');\nawait context.sync();\n\nconst oldLabel = 'Resolution for Application to Proceed to the ALC';\nconst newLabel = 'What is the outcome of the Board/Council resolution?';\nif (resolutionPara.text === oldLabel) {\n  resolutionPara.getRange().insertText(newLabel, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The L/FNG Review template has several table-label cells whose paragraph\n# mark run properties (w:pPr/w:rPr) carry an accidentally duplicated\n# <w:b/> element (bold turned on twice). We re-apply bold on the range of\n# each of those label cells so the document model re-serializes the\n# (single, deduplicated) bold flag. We also update the \"Resolution for\n# Application to Proceed to the ALC\" label text to the new question\n# wording, per the commit \"Update L/FNG Review Form Labels\".\n\n$d = $word.ActiveDocument\n\n# Tables() / Cell() are 1-based: [tableIndex, rowIndex, columnIndex] for\n# every label cell that has the duplicated bold flag in its paragraph\n# mark run properties.\n$targets = @(\n    @(1, 1, 1), # Application ID:\n    @(1, 2, 1), # Application Type:\n    @(1, 3, 1), # Status:\n    @(1, 4, 1), # Applicant:\n    @(1, 5, 1), # Local/ First Nation Government:\n    @(2, 1, 1), # Local Government File Number\n    @(2, 2, 1), # First Name\n    @(2, 3, 1), # Last Name\n    @(2, 4, 1), # Position\n    @(2, 5, 1), # Department\n    @(2, 6, 1), # Phone Number\n    @(2, 7, 1), # Email\n    @(3, 1, 1), # Resolution for Application to Proceed to the ALC\n    @(4, 1, 1), # Type (visible attachments table header)\n    @(4, 1, 2), # File Name (visible attachments table header)\n    @(5, 1, 1), # Type (no-data attachments table header)\n    @(5, 1, 2)  # File Name (no-data attachments table header)\n)\n\nforeach ($target in $targets) {\n    $tableIndex = $target[0]\n    $rowIndex = $target[1]\n    $columnIndex = $target[2]\n    $cell = $d.Tables.Item($tableIndex).Cell($rowIndex, $columnIndex)\n    # Re-asserting bold forces the cell's paragraph formatting (w:pPr/w:rPr\n    # and w:r/w:rPr) to be rewritten without the stray duplicate <w:b/>.\n    $cell.Range.Font.Bold = 1\n}\n\n# Update the \"Resolution\" row label text, after confirming we have the\n# right cell.\n$oldLabel = \"Resolution for Application to Proceed to the ALC\"\n$newLabel = \"What is the outcome of the Board/Council resolution?\"\n$resolutionCell = $d.Tables.Item(3).Cell(1, 1)\n$resolutionText = $resolutionCell.Range.Text.TrimEnd([char]13, [char]7)\nif ($resolutionText -eq $oldLabel) {\n    $resolutionCell.Range.Text = $newLabel\n}\n"}
